$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape on slide 1 that holds the "Hello Salhi Karem" line
# (the "ZoneTexte 12" textbox with the "Realise par :" / name paragraph).
foreach ($shape in $s.Shapes) {
    if ($shape.HasTextFrame) {
        $tf = $shape.TextFrame
        $tr = $tf.TextRange
        if ($tr.Text -like "*Hello Salhi*") {
            # The shape auto-fits to its text (a:spAutoFit); remember its
            # current height so we can restore it after the edit (the source
            # edit did not change the shape's stored extent).
            $origHeight = $shape.Height

            # Find the paragraph that contains "Hello Salhi"
            $paraCount = $tr.Paragraphs().Count
            for ($i = 1; $i -le $paraCount; $i++) {
                $para = $tr.Paragraphs($i)
                if ($para.Text -like "*Hello Salhi*") {
                    # The paragraph currently reads "Hello Salhi Karem    ".
                    # Replace the leading "Hello" (chars 1-5) with "BY", then
                    # re-type the separating space and "Salhi" so they land
                    # in their own runs, matching: "BY" + " " + "Salhi".
                    $c1 = $para.Characters(1, 5)
                    $c1.Text = "BY"

                    $c2 = $para.Characters(3, 1)
                    $c2.Text = " "

                    $c3 = $para.Characters(4, 5)
                    $c3.Text = "Salhi"
                }
            }

            # Restore the shape's height (the auto-fit box shrank because
            # the replacement text is shorter than the original).
            $shape.Height = $origHeight
        }
    }
}
